# Ccl12-Ccr1 LR-pair sheet: refresh rows 2-3 with recomputed NATMI values
# (per Dr Hou's advice) and add the remaining 7 sending/target cluster
# combinations (rows 4-10) to complete the full 3x3 matrix across the
# FAPs / M2 / sCs sending clusters and ECs / M2 / sCs target clusters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Ccl12"
$ws.Range("C2").Value = "Ccr1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.857148
$ws.Range("H2").Value = 2.571444
$ws.Range("I2").Value = 0.04787301688248034
$ws.Range("J2").Value = 0.04787301688248034
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.006600000000000001
$ws.Range("N2").Value = 0.0198
$ws.Range("O2").Value = 0.001099924410750217
$ws.Range("P2").Value = 0.001099924410750217
$ws.Range("Q2").Value = 0.005657176800000001
$ws.Range("R2").Value = 0.0509145912
$ws.Range("S2").Value = 0.00005265669988529737
$ws.Range("T2").Value = 0.00005265669988529737

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Ccl12"
$ws.Range("C3").Value = "Ccr1"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.857148
$ws.Range("H3").Value = 2.571444
$ws.Range("I3").Value = 0.04787301688248034
$ws.Range("J3").Value = 0.04787301688248034
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.926252666666667
$ws.Range("N3").Value = 17.778758
$ws.Range("O3").Value = 0.9876409049000355
$ws.Range("P3").Value = 0.9876409049000355
$ws.Range("Q3").Value = 5.079675620728
$ws.Range("R3").Value = 45.717080586552
$ws.Range("S3").Value = 0.04728134971410756
$ws.Range("T3").Value = 0.04728134971410757

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Ccl12"
$ws.Range("C4").Value = "Ccr1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.857148
$ws.Range("H4").Value = 2.571444
$ws.Range("I4").Value = 0.04787301688248034
$ws.Range("J4").Value = 0.04787301688248034
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.06755966666666667
$ws.Range("N4").Value = 0.202679
$ws.Range("O4").Value = 0.0112591706892143
$ws.Range("P4").Value = 0.0112591706892143
$ws.Range("Q4").Value = 0.057908633164
$ws.Range("R4").Value = 0.521177698476
$ws.Range("S4").Value = 0.0005390104684874841
$ws.Range("T4").Value = 0.0005390104684874842

# Row 5
$ws.Range("A5").Value = "M2"
$ws.Range("B5").Value = "Ccl12"
$ws.Range("C5").Value = "Ccr1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 15.94278666666667
$ws.Range("H5").Value = 47.82836
$ws.Range("I5").Value = 0.8904288352152905
$ws.Range("J5").Value = 0.8904288352152906
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.006600000000000001
$ws.Range("N5").Value = 0.0198
$ws.Range("O5").Value = 0.001099924410750217
$ws.Range("P5").Value = 0.001099924410750217
$ws.Range("Q5").Value = 0.105222392
$ws.Range("R5").Value = 0.9470015280000001
$ws.Range("S5").Value = 0.0009794044118891804
$ws.Range("T5").Value = 0.0009794044118891804

# Row 6
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Ccl12"
$ws.Range("C6").Value = "Ccr1"
$ws.Range("D6").Value = "M2"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 15.94278666666667
$ws.Range("H6").Value = 47.82836
$ws.Range("I6").Value = 0.8904288352152905
$ws.Range("J6").Value = 0.8904288352152906
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 5.926252666666667
$ws.Range("N6").Value = 17.778758
$ws.Range("O6").Value = 0.9876409049000355
$ws.Range("P6").Value = 0.9876409049000355
$ws.Range("Q6").Value = 94.48098199743112
$ws.Range("R6").Value = 850.3288379768801
$ws.Range("S6").Value = 0.8794239405611142
$ws.Range("T6").Value = 0.8794239405611143

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Ccl12"
$ws.Range("C7").Value = "Ccr1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 15.94278666666667
$ws.Range("H7").Value = 47.82836
$ws.Range("I7").Value = 0.8904288352152905
$ws.Range("J7").Value = 0.8904288352152906
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.06755966666666667
$ws.Range("N7").Value = 0.202679
$ws.Range("O7").Value = 0.0112591706892143
$ws.Range("P7").Value = 0.0112591706892143
$ws.Range("Q7").Value = 1.077089352937778
$ws.Range("R7").Value = 9.69380417644
$ws.Range("S7").Value = 0.01002549024228723
$ws.Range("T7").Value = 0.01002549024228723

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Ccl12"
$ws.Range("C8").Value = "Ccr1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 1.104681666666667
$ws.Range("H8").Value = 3.314045
$ws.Range("I8").Value = 0.06169814790222908
$ws.Range("J8").Value = 0.06169814790222908
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.006600000000000001
$ws.Range("N8").Value = 0.0198
$ws.Range("O8").Value = 0.001099924410750217
$ws.Range("P8").Value = 0.001099924410750217
$ws.Range("Q8").Value = 0.007290899000000002
$ws.Range("R8").Value = 0.065618091
$ws.Range("S8").Value = 0.00006786329897573906
$ws.Range("T8").Value = 0.00006786329897573905

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Ccl12"
$ws.Range("C9").Value = "Ccr1"
$ws.Range("D9").Value = "M2"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 1.104681666666667
$ws.Range("H9").Value = 3.314045
$ws.Range("I9").Value = 0.06169814790222908
$ws.Range("J9").Value = 0.06169814790222908
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 5.926252666666667
$ws.Range("N9").Value = 17.778758
$ws.Range("O9").Value = 0.9876409049000355
$ws.Range("P9").Value = 0.9876409049000355
$ws.Range("Q9").Value = 6.546622672901112
$ws.Range("R9").Value = 58.91960405611
$ws.Range("S9").Value = 0.06093561462481376
$ws.Range("T9").Value = 0.06093561462481376

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Ccl12"
$ws.Range("C10").Value = "Ccr1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 1.104681666666667
$ws.Range("H10").Value = 3.314045
$ws.Range("I10").Value = 0.06169814790222908
$ws.Range("J10").Value = 0.06169814790222908
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.06755966666666667
$ws.Range("N10").Value = 0.202679
$ws.Range("O10").Value = 0.0112591706892143
$ws.Range("P10").Value = 0.0112591706892143
$ws.Range("Q10").Value = 0.0746319251727778
$ws.Range("R10").Value = 0.6716873265550001
$ws.Range("S10").Value = 0.0006946699784395866
$ws.Range("T10").Value = 0.0006946699784395866

